$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 39, pushing existing rows
# 39-50 down to 41-52 (dimension grows from A1:R50 to A1:R52).
$ws.Rows("39:40").Insert()

# --- New row 39 ---
$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44559
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 100112030
$ws.Range("G39").Value = "Poroto granado"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 360
$ws.Range("K39").Value = 30000
$ws.Range("L39").Value = 31000
$ws.Range("M39").Value = 30500
$ws.Range("N39").Value = "$/malla 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 1220
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"

# --- New row 40 ---
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "Terminal La Palmera de La Serena"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44559
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 100112030
$ws.Range("G40").Value = "Poroto granado"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 32000
$ws.Range("L40").Value = 33000
$ws.Range("M40").Value = 32500
$ws.Range("N40").Value = "$/malla 25 kilos"
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 1300
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
